$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain text (matching original inlineStr text cells) while updating values
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.49%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "18"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.93%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "18"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.696"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.09%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "18"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06080"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.64%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "18"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.737"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.77%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "18"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8506"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.07%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "18"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9054"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.06%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "18"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1411"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.28%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "18"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05030"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.99%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "18"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07100"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.17%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "18"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03153"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.06%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "18"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09018"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.13%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "18"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001532"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.56%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "18"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006073"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.43%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "18"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006078"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.19%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "18"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.455"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.10%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "18"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.171"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.19%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "18"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.197"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.31%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "18"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.38%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "18"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.092"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.22%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "18"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04255"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.48%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "18"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001180"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "18"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004136"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.78%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.08%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "18"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "6.63%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "18"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "18"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "18"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "18"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "18"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "18"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "18"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "18"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "18"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "18"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "18"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "18"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "18"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03918"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.44%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "18"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.19%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "18"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004175"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.76%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "18"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "18"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01148"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-17.39%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "18"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005103"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.72%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "18"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "18"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "18"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1653"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.80%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "18"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "18"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "18"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "18"
